$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 135.08333
$ws.Range("H6").Value = 12000
$ws.Range("I6").Value = 12000
$ws.Range("K6").Value = 36000
$ws.Range("M6").Value = -35888
$ws.Range("H11").Value = 110550.57
$ws.Range("I11").Value = 110550.57
$ws.Range("K11").Value = 110550.57
$ws.Range("M11").Value = -110410.57
$ws.Range("H19").Value = 1543.0476
$ws.Range("I19").Value = 867.2
$ws.Range("K19").Value = 867.2
$ws.Range("M19").Value = -692.2
$ws.Range("H33").Value = 148
$ws.Range("I33").Value = 152.45454
$ws.Range("J33").Value = 99
$ws.Range("K33").Value = 152.45454
$ws.Range("L33").Value = 99
$ws.Range("M33").Value = 76.54545999999999
$ws.Range("N33").Value = -557
$ws.Range("H43").Value = 3782.8462
$ws.Range("I43").Value = 3847.75
$ws.Range("J43").Value = 3754
$ws.Range("K43").Value = 3847.75
$ws.Range("L43").Value = 3754
$ws.Range("M43").Value = -3778.75
$ws.Range("N43").Value = -3892
$ws.Range("H62").Value = 3455.2
$ws.Range("I62").Value = 2750
$ws.Range("K62").Value = 2750
$ws.Range("M62").Value = -2126
$ws.Range("H65").Value = 3455.2
$ws.Range("I65").Value = 2750
$ws.Range("K65").Value = 13750
$ws.Range("M65").Value = -10630
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H92").Value = 1385.0588
$ws.Range("I92").Value = 2106.8
$ws.Range("K92").Value = 2106.8
$ws.Range("M92").Value = -858.8000000000002
$ws.Range("H97").Value = 3000
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -8504
$ws.Range("H98").Value = 1358.5454
$ws.Range("I98").Value = 1358.5454
$ws.Range("K98").Value = 1358.5454
$ws.Range("M98").Value = 139.4546
$ws.Range("H100").Value = 1681.25
$ws.Range("I100").Value = 1575
$ws.Range("K100").Value = 1575
$ws.Range("M100").Value = -1034
$ws.Range("H122").Value = 1358.5454
$ws.Range("I122").Value = 1358.5454
$ws.Range("K122").Value = 4075.6362
$ws.Range("M122").Value = -1625.6362
$ws.Range("H125").Value = 1891.409
$ws.Range("I125").Value = 1855.7333
$ws.Range("K125").Value = 16701.5997
$ws.Range("M125").Value = -14241.5997
$ws.Range("H127").Value = 80357
$ws.Range("I127").Value = 80357
$ws.Range("K127").Value = 241071
$ws.Range("M127").Value = -236111
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 21887.455
$ws.Range("I132").Value = 21887.455
$ws.Range("K132").Value = 65662.36500000001
$ws.Range("M132").Value = -63132.36500000001
$ws.Range("H133").Value = 89779.5
$ws.Range("I133").Value = 89779
$ws.Range("J133").Value = 89780
$ws.Range("K133").Value = 89779
$ws.Range("L133").Value = 89780
$ws.Range("M133").Value = -84719
$ws.Range("N133").Value = -99900
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("H137").Value = 5891265.5
$ws.Range("J137").Value = 12344.429
$ws.Range("L137").Value = 37033.287
$ws.Range("N137").Value = -42133.287
$ws.Range("H138").Value = 6596.32
$ws.Range("J138").Value = 7034.9414
$ws.Range("L138").Value = 21104.8242
$ws.Range("N138").Value = -31384.8242
$ws.Range("H139").Value = 78776.625
$ws.Range("J139").Value = 78776.625
$ws.Range("L139").Value = 78776.625
$ws.Range("N139").Value = -89056.625
$ws.Range("H141").Value = 6646.9473
$ws.Range("I141").Value = 3130.5
$ws.Range("J141").Value = 10554.111
$ws.Range("K141").Value = 9391.5
$ws.Range("L141").Value = 31662.333
$ws.Range("M141").Value = -4211.5
$ws.Range("N141").Value = -42022.333
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("N136").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 880171
$ws.Range("I32").Value = 894102.6
$ws.Range("K32").Value = 894102.6
$ws.Range("M32").Value = -893815.6
$ws.Range("H45").Value = 3561.2727
$ws.Range("I45").Value = 3396.5715
$ws.Range("K45").Value = 3396.5715
$ws.Range("M45").Value = -3019.5715
$ws.Range("H61").Value = 3705824.2
$ws.Range("I61").Value = 2779620.8
$ws.Range("K61").Value = 2779620.8
$ws.Range("M61").Value = -2779408.8
$ws.Range("H74").Value = 1618610.5
$ws.Range("I74").Value = 2294530.8
$ws.Range("K74").Value = 2294530.8
$ws.Range("M74").Value = -2293656.8
$ws.Range("H77").Value = 1618610.5
$ws.Range("I77").Value = 2294530.8
$ws.Range("K77").Value = 11472654
$ws.Range("M77").Value = -11468286
$ws.Range("H125").Value = 38990
$ws.Range("J125").Value = 38990
$ws.Range("L125").Value = 38990
$ws.Range("N125").Value = -48830
$ws.Range("H132").Value = 4190.4287
$ws.Range("I132").Value = 2948.8928
$ws.Range("J132").Value = 5845.8096
$ws.Range("K132").Value = 8846.678400000001
$ws.Range("L132").Value = 17537.4288
$ws.Range("M132").Value = -6316.678400000001
$ws.Range("N132").Value = -22597.4288
$ws.Range("H136").Value = 3705824.2
$ws.Range("I136").Value = 2779620.8
$ws.Range("K136").Value = 8338862.399999999
$ws.Range("M136").Value = -8336312.399999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 269474.28
$ws.Range("I20").Value = 337644.03
$ws.Range("K20").Value = 337644.03
$ws.Range("M20").Value = -337397.03
$ws.Range("H82").Value = 20617.223
$ws.Range("J82").Value = 51500
$ws.Range("L82").Value = 51500
$ws.Range("N82").Value = -52266
$ws.Range("H85").Value = 20617.223
$ws.Range("J85").Value = 51500
$ws.Range("L85").Value = 51500
$ws.Range("N85").Value = -54152
$ws.Range("H134").Value = 4112980
$ws.Range("I134").Value = 3482701.8
$ws.Range("J134").Value = 5562619.5
$ws.Range("K134").Value = 10448105.4
$ws.Range("L134").Value = 16687858.5
$ws.Range("M134").Value = -10445570.4
$ws.Range("N134").Value = -16692928.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 64900.688
$ws.Range("I16").Value = 2560.7334
$ws.Range("J16").Value = 1000000
$ws.Range("K16").Value = 2560.7334
$ws.Range("L16").Value = 1000000
$ws.Range("M16").Value = -2273.7334
$ws.Range("N16").Value = -1000574
$ws.Range("H31").Value = 473581.1
$ws.Range("I31").Value = 1016222.1
$ws.Range("J31").Value = 3292.2334
$ws.Range("K31").Value = 1016222.1
$ws.Range("L31").Value = 3292.2334
$ws.Range("M31").Value = -1015927.1
$ws.Range("N31").Value = -3882.2334
$ws.Range("H34").Value = 473581.1
$ws.Range("I34").Value = 1016222.1
$ws.Range("J34").Value = 3292.2334
$ws.Range("K34").Value = 1016222.1
$ws.Range("L34").Value = 3292.2334
$ws.Range("M34").Value = -1016020.1
$ws.Range("N34").Value = -3696.2334
$ws.Range("H58").Value = 9140471
$ws.Range("J58").Value = 2176958.5
$ws.Range("L58").Value = 2176958.5
$ws.Range("N58").Value = -2177364.5
$ws.Range("H86").Value = 40365.44
$ws.Range("J86").Value = 9492.904
$ws.Range("L86").Value = 9492.904
$ws.Range("N86").Value = -11738.904
$ws.Range("H89").Value = 40365.44
$ws.Range("J89").Value = 9492.904
$ws.Range("L89").Value = 47464.52
$ws.Range("N89").Value = -58696.52
$ws.Range("H93").Value = 17599
$ws.Range("I93").Value = 17599
$ws.Range("K93").Value = 17599
$ws.Range("M93").Value = -15727
$ws.Range("H94").Value = 21281.4
$ws.Range("J94").Value = 1802.6666
$ws.Range("L94").Value = 1802.6666
$ws.Range("N94").Value = -2704.6666
$ws.Range("H99").Value = 824065.9
$ws.Range("I99").Value = 68553.336
$ws.Range("K99").Value = 68553.336
$ws.Range("M99").Value = -67055.336
$ws.Range("H113").Value = 64900.688
$ws.Range("I113").Value = 2560.7334
$ws.Range("J113").Value = 1000000
$ws.Range("K113").Value = 2560.7334
$ws.Range("L113").Value = 1000000
$ws.Range("M113").Value = -390.7334000000001
$ws.Range("N113").Value = -1004340
$ws.Range("H117").Value = 50000
$ws.Range("I117").Value = 50000
$ws.Range("K117").Value = 50000
$ws.Range("M117").Value = -45411
$ws.Range("H122").Value = 19566.545
$ws.Range("I122").Value = 512.8
$ws.Range("J122").Value = 35444.668
$ws.Range("K122").Value = 1538.4
$ws.Range("L122").Value = 106334.004
$ws.Range("M122").Value = 911.6000000000001
$ws.Range("N122").Value = -111234.004
$ws.Range("H126").Value = 824065.9
$ws.Range("I126").Value = 68553.336
$ws.Range("K126").Value = 205660.008
$ws.Range("M126").Value = -203190.008
$ws.Range("H132").Value = 2551.4285
$ws.Range("I132").Value = 2460.1667
$ws.Range("K132").Value = 7380.500100000001
$ws.Range("M132").Value = -4850.500100000001
$ws.Range("H134").Value = 5449.8296
$ws.Range("I134").Value = 4165.385
$ws.Range("J134").Value = 5940.9414
$ws.Range("K134").Value = 12496.155
$ws.Range("L134").Value = 17822.8242
$ws.Range("M134").Value = -9961.155000000001
$ws.Range("N134").Value = -22892.8242
$ws.Range("H136").Value = 9140471
$ws.Range("J136").Value = 2176958.5
$ws.Range("L136").Value = 6530875.5
$ws.Range("N136").Value = -6535975.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 97.8421
$ws.Range("I2").Value = 20
$ws.Range("K2").Value = 120
$ws.Range("M2").Value = -7
$ws.Range("H4").Value = 12226028
$ws.Range("I4").Value = 17357468
$ws.Range("J4").Value = 252666.5
$ws.Range("K4").Value = 52072404
$ws.Range("L4").Value = 757999.5
$ws.Range("M4").Value = -52072292
$ws.Range("N4").Value = -758223.5
$ws.Range("H5").Value = 1282585.8
$ws.Range("I5").Value = 1374134.6
$ws.Range("J5").Value = 1163572.2
$ws.Range("K5").Value = 4122403.8
$ws.Range("L5").Value = 3490716.6
$ws.Range("M5").Value = -4122291.8
$ws.Range("N5").Value = -3490940.6
$ws.Range("H18").Value = 77570.766
$ws.Range("I18").Value = 143055
$ws.Range("K18").Value = 429165
$ws.Range("M18").Value = -428996
$ws.Range("H32").Value = 71751576
$ws.Range("I32").Value = 138500660
$ws.Range("K32").Value = 415501980
$ws.Range("M32").Value = -415501697
$ws.Range("H35").Value = 111111160
$ws.Range("I35").Value = 72.5
$ws.Range("J35").Value = 333333340
$ws.Range("K35").Value = 217.5
$ws.Range("L35").Value = 1000000020
$ws.Range("M35").Value = 70.5
$ws.Range("N35").Value = -1000000596
$ws.Range("H36").Value = 1310.2
$ws.Range("I36").Value = 183.66667
$ws.Range("J36").Value = 3000
$ws.Range("K36").Value = 551.00001
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = -382.00001
$ws.Range("N36").Value = -9338
$ws.Range("H41").Value = 275
$ws.Range("I41").Value = 300
$ws.Range("J41").Value = 250
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 750
$ws.Range("M41").Value = -562
$ws.Range("N41").Value = -1426
$ws.Range("H68").Value = 4227.9375
$ws.Range("J68").Value = 4355.6333
$ws.Range("L68").Value = 13066.8999
$ws.Range("N68").Value = -14688.8999
$ws.Range("H71").Value = 4227.9375
$ws.Range("J71").Value = 4355.6333
$ws.Range("L71").Value = 39200.6997
$ws.Range("N71").Value = -47312.6997
$ws.Range("H93").Value = 175.33333
$ws.Range("I93").Value = 175.33333
$ws.Range("K93").Value = 525.99999
$ws.Range("M93").Value = 1346.00001
$ws.Range("H135").Value = 1282585.8
$ws.Range("I135").Value = 1374134.6
$ws.Range("J135").Value = 1163572.2
$ws.Range("K135").Value = 12367211.4
$ws.Range("L135").Value = 10472149.8
$ws.Range("M135").Value = -12364676.4
$ws.Range("N135").Value = -10477219.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 18000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("H70").Value = 14138.381
$ws.Range("I70").Value = 30701.857
$ws.Range("K70").Value = 30701.857
$ws.Range("M70").Value = -30431.857
$ws.Range("H73").Value = 14138.381
$ws.Range("I73").Value = 30701.857
$ws.Range("K73").Value = 30701.857
$ws.Range("M73").Value = -29765.857
$ws.Range("H80").Value = 1370.8334
$ws.Range("I80").Value = 1249.3334
$ws.Range("K80").Value = 1249.3334
$ws.Range("M80").Value = -251.3334
$ws.Range("H83").Value = 1370.8334
$ws.Range("I83").Value = 1249.3334
$ws.Range("K83").Value = 6246.666999999999
$ws.Range("M83").Value = -1254.666999999999
$ws.Range("H107").Value = 1013.4286
$ws.Range("I107").Value = 961
$ws.Range("J107").Value = 1118.2858
$ws.Range("K107").Value = 961
$ws.Range("L107").Value = 1118.2858
$ws.Range("M107").Value = 959
$ws.Range("N107").Value = -4958.2858
$ws.Range("H113").Value = 5167.5
$ws.Range("I113").Value = 5081
$ws.Range("J113").Value = 5600
$ws.Range("K113").Value = 5081
$ws.Range("L113").Value = 5600
$ws.Range("M113").Value = -2911
$ws.Range("N113").Value = -9940
$ws.Range("H122").Value = 10148113
$ws.Range("J122").Value = 1947.75
$ws.Range("L122").Value = 5843.25
$ws.Range("N122").Value = -10743.25
$ws.Range("H132").Value = 14059.948
$ws.Range("I132").Value = 12083.63
$ws.Range("J132").Value = 18506.666
$ws.Range("K132").Value = 36250.89
$ws.Range("L132").Value = 55519.99800000001
$ws.Range("M132").Value = -33720.89
$ws.Range("N132").Value = -60579.99800000001
$ws.Range("H139").Value = 102321
$ws.Range("J139").Value = 102321
$ws.Range("L139").Value = 102321
$ws.Range("N139").Value = -112601
$ws.Range("M47").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2799.6667
$ws.Range("I7").Value = 2799.6667
$ws.Range("K7").Value = 2799.6667
$ws.Range("M7").Value = -2687.6667
$ws.Range("H16").Value = 1522
$ws.Range("J16").Value = 1563.25
$ws.Range("L16").Value = 1563.25
$ws.Range("N16").Value = -1903.25
$ws.Range("H22").Value = 2699.75
$ws.Range("I22").Value = 1900
$ws.Range("J22").Value = 3099.625
$ws.Range("K22").Value = 1900
$ws.Range("L22").Value = 3099.625
$ws.Range("M22").Value = -1605
$ws.Range("N22").Value = -3689.625
$ws.Range("H27").Value = 2699.75
$ws.Range("I27").Value = 1900
$ws.Range("J27").Value = 3099.625
$ws.Range("K27").Value = 1900
$ws.Range("L27").Value = 3099.625
$ws.Range("M27").Value = -1793
$ws.Range("N27").Value = -3313.625
$ws.Range("H82").Value = 1737.6428
$ws.Range("I82").Value = 1673.2
$ws.Range("J82").Value = 1898.75
$ws.Range("K82").Value = 1673.2
$ws.Range("L82").Value = 1898.75
$ws.Range("M82").Value = -1312.2
$ws.Range("N82").Value = -2620.75
$ws.Range("H85").Value = 1737.6428
$ws.Range("I85").Value = 1673.2
$ws.Range("J85").Value = 1898.75
$ws.Range("K85").Value = 1673.2
$ws.Range("L85").Value = 1898.75
$ws.Range("M85").Value = -425.2
$ws.Range("N85").Value = -4394.75
$ws.Range("H100").Value = 2142.1875
$ws.Range("J100").Value = 2533.3333
$ws.Range("L100").Value = 2533.3333
$ws.Range("N100").Value = -3615.3333
$ws.Range("H126").Value = 2799.6667
$ws.Range("I126").Value = 2799.6667
$ws.Range("K126").Value = 8399.000100000001
$ws.Range("M126").Value = -5929.000100000001
$ws.Range("H132").Value = 2979597.5
$ws.Range("I132").Value = 4169631.8
$ws.Range("J132").Value = 4511.875
$ws.Range("K132").Value = 12508895.4
$ws.Range("L132").Value = 13535.625
$ws.Range("M132").Value = -12506365.4
$ws.Range("N132").Value = -18595.625
$ws.Range("H136").Value = 93500584
$ws.Range("I136").Value = 83000780
$ws.Range("K136").Value = 249002340
$ws.Range("M136").Value = -248999790

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 151323.28
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 151323.28
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 302646.56
$ws.Range("N81").Value = -304768.56
$ws.Range("H84").Value = 151323.28
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 151323.28
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 1513232.8
$ws.Range("N84").Value = -1523840.8
$ws.Range("H107").Value = 1427.08
$ws.Range("I107").Value = 761.5
$ws.Range("J107").Value = 6308
$ws.Range("K107").Value = 2284.5
$ws.Range("L107").Value = 18924
$ws.Range("M107").Value = -364.5
$ws.Range("N107").Value = -22764
$ws.Range("H113").Value = 1799.6
$ws.Range("I113").Value = 1244.6666
$ws.Range("J113").Value = 2253.6365
$ws.Range("K113").Value = 3733.9998
$ws.Range("L113").Value = 6760.9095
$ws.Range("M113").Value = -1563.9998
$ws.Range("N113").Value = -11100.9095
$ws.Range("H126").Value = 2994.3
$ws.Range("I126").Value = 3024.25
$ws.Range("K126").Value = 9072.75
$ws.Range("M126").Value = -6602.75
$ws.Range("H132").Value = 7578001.5
$ws.Range("I132").Value = 8335459
$ws.Range("K132").Value = 25006377
$ws.Range("M132").Value = -25003847
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()
